$d = $word.ActiveDocument

# Locate the paragraph beginning "Manager begins the script with a step index..."
$startRange = $d.Content.Duplicate
$found1 = $startRange.Find.Execute("Manager begins the script with a step index", $false, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find start paragraph"
}
$startPara = $startRange.Paragraphs(1)

# Locate the paragraph ending "...connected to the corresponding slots in Manager."
$endRange = $d.Content.Duplicate
$found2 = $endRange.Find.Execute("connected to the corresponding slots in Manager.", $false, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find end paragraph"
}
$endPara = $endRange.Paragraphs(1)

# Remove all four list-item paragraphs (including their paragraph marks) as a single range.
$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete()
